$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Propagate the existing date number-format (currently s="1" on
# C2:C4) onto the new, empty D2:D4 cells before the C column content
# is overwritten below.
# ------------------------------------------------------------------
$ws.Range("C2").Copy() | Out-Null
$ws.Range("D2").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Header row
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Internal Name"
$ws.Range("B1").Value = "External Name"
$ws.Range("C1").Value = "Industry"

# ------------------------------------------------------------------
# Data rows
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Energy Company"
$ws.Range("B2").Value = "Power Inc."
$ws.Range("C2").Value = "Energy"
$ws.Range("C2").Style = "Normal"

$ws.Range("A3").Value = "Engineering Firm"
$ws.Range("B3").Value = "Acme Inc."
$ws.Range("C3").Value = "Engineering"
$ws.Range("C3").Style = "Normal"

$ws.Range("A4").Value = "RPA Provider"
$ws.Range("B4").Value = "UiPath"
$ws.Range("C4").Value = "Engineering"
$ws.Range("C4").Style = "Normal"

# ------------------------------------------------------------------
# Column widths (approximating the column widths recorded in the
# target workbook; Excel snaps column widths to whole-pixel
# increments, so these are the closest achievable values).
# ------------------------------------------------------------------
$ws.Range("A1").ColumnWidth = 15.33
$ws.Range("B1").ColumnWidth = 17.83
$ws.Range("C1").ColumnWidth = 16.0

Write-Host "edits applied"
